$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution timestamp recorded for the existing row (A2) --
# the scheduled task re-ran and produced a marginally different
# fractional-seconds value for the same logged moment.
$ws.Range("A2").Value = 45873.33357223379

# Append the new measurement row captured by the latest scheduled run.
$ws.Range("A3").Value = 45873.37521892216
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 15.56
$ws.Range("E3").Value = 87.52
$ws.Range("F3").Value = 333.98
$ws.Range("G3").Value = 8.73
$ws.Range("H3").Value = "SE"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "09:00:18"
